$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values reflect a re-run of the NATMI analysis with refreshed TPM
# expression inputs. Only the numeric measurement columns (G..T) change;
# identifiers in columns A-D and the raw counts in columns E, F, K, L are
# unaffected.

$ws.Range("G2").Value2 = 73.823705
$ws.Range("H2").Value2 = 221.471115
$ws.Range("I2").Value2 = 0.260518167190625
$ws.Range("J2").Value2 = 0.260518167190625
$ws.Range("M2").Value2 = 91.51130433333333
$ws.Range("N2").Value2 = 274.533913
$ws.Range("O2").Value2 = 0.9685519820468944
$ws.Range("P2").Value2 = 0.9685519820468945
$ws.Range("Q2").Value2 = 6755.703535269222
$ws.Range("R2").Value2 = 60801.33181742299
$ws.Range("S2").Value2 = 0.2523253871917041
$ws.Range("T2").Value2 = 0.2523253871917041
$ws.Range("G3").Value2 = 73.823705
$ws.Range("H3").Value2 = 221.471115
$ws.Range("I3").Value2 = 0.260518167190625
$ws.Range("J3").Value2 = 0.260518167190625
$ws.Range("O3").Value2 = 0.001425786415744213
$ws.Range("P3").Value2 = 0.001425786415744214
$ws.Range("Q3").Value2 = 9.944938947960001
$ws.Range("R3").Value2 = 89.50445053164
$ws.Range("S3").Value2 = 0.000371443263834973
$ws.Range("T3").Value2 = 0.0003714432638349731
$ws.Range("G4").Value2 = 73.823705
$ws.Range("H4").Value2 = 221.471115
$ws.Range("I4").Value2 = 0.260518167190625
$ws.Range("J4").Value2 = 0.260518167190625
$ws.Range("M4").Value2 = 2.836578333333333
$ws.Range("N4").Value2 = 8.509734999999999
$ws.Range("O4").Value2 = 0.03002223153736139
$ws.Range("P4").Value2 = 0.03002223153736139
$ws.Range("Q4").Value2 = 209.4067220893917
$ws.Range("R4").Value2 = 1884.660498804525
$ws.Range("S4").Value2 = 0.007821336735085971
$ws.Range("T4").Value2 = 0.007821336735085971
$ws.Range("I5").Value2 = 0.5089642665472768
$ws.Range("J5").Value2 = 0.5089642665472768
$ws.Range("M5").Value2 = 91.51130433333333
$ws.Range("N5").Value2 = 274.533913
$ws.Range("O5").Value2 = 0.9685519820468944
$ws.Range("P5").Value2 = 0.9685519820468945
$ws.Range("Q5").Value2 = 13198.35669012367
$ws.Range("R5").Value2 = 118785.2102111131
$ws.Range("S5").Value2 = 0.4929583491554088
$ws.Range("T5").Value2 = 0.4929583491554089
$ws.Range("I6").Value2 = 0.5089642665472768
$ws.Range("J6").Value2 = 0.5089642665472768
$ws.Range("O6").Value2 = 0.001425786415744213
$ws.Range("P6").Value2 = 0.001425786415744214
$ws.Range("S6").Value2 = 0.0007256743373423243
$ws.Range("T6").Value2 = 0.0007256743373423244
$ws.Range("I7").Value2 = 0.5089642665472768
$ws.Range("J7").Value2 = 0.5089642665472768
$ws.Range("M7").Value2 = 2.836578333333333
$ws.Range("N7").Value2 = 8.509734999999999
$ws.Range("O7").Value2 = 0.03002223153736139
$ws.Range("P7").Value2 = 0.03002223153736139
$ws.Range("S7").Value2 = 0.01528024305452566
$ws.Range("T7").Value2 = 0.01528024305452566
$ws.Range("G8").Value2 = 65.32235733333334
$ws.Range("I8").Value2 = 0.2305175662620982
$ws.Range("J8").Value2 = 0.2305175662620982
$ws.Range("M8").Value2 = 91.51130433333333
$ws.Range("N8").Value2 = 274.533913
$ws.Range("O8").Value2 = 0.9685519820468944
$ws.Range("P8").Value2 = 0.9685519820468945
$ws.Range("Q8").Value2 = 5977.734121701415
$ws.Range("R8").Value2 = 53799.60709531274
$ws.Range("S8").Value2 = 0.2232682456997815
$ws.Range("T8").Value2 = 0.2232682456997815
$ws.Range("G9").Value2 = 65.32235733333334
$ws.Range("I9").Value2 = 0.2305175662620982
$ws.Range("J9").Value2 = 0.2305175662620982
$ws.Range("O9").Value2 = 0.001425786415744213
$ws.Range("P9").Value2 = 0.001425786415744214
$ws.Range("Q9").Value2 = 8.799705401088001
$ws.Range("R9").Value2 = 79.19734860979202
$ws.Range("S9").Value2 = 0.0003286688145669161
$ws.Range("T9").Value2 = 0.0003286688145669163
$ws.Range("G10").Value2 = 65.32235733333334
$ws.Range("I10").Value2 = 0.2305175662620982
$ws.Range("J10").Value2 = 0.2305175662620982
$ws.Range("M10").Value2 = 2.836578333333333
$ws.Range("N10").Value2 = 8.509734999999999
$ws.Range("O10").Value2 = 0.03002223153736139
$ws.Range("P10").Value2 = 0.03002223153736139
$ws.Range("S10").Value2 = 0.006920651747749759
$ws.Range("T10").Value2 = 0.00692065174774976
